$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace accented Portuguese text with LaTeX-escaped equivalents in the
# "Institution" column (column E).
$ws.Range("E2").Value = "Funda\c{c}\~ao Calouste Gulbenkian"
$ws.Range("E3").Value = "Junta Nacional de Investiga\c{c}\~ao Cientifica e Tecnologica "
$ws.Range("E4").Value = "Junta Nacional de Investiga\c{c}\~ao Cientifica e Tecnologica "
$ws.Range("E6").Value = "Funda\c{c}\~ao para a Ci\^encia e Tecnologia"

# Move the active selection from E6 to E7, matching the final cursor
# position recorded in the saved workbook.
$ws.Range("E7").Select()
